{"js": "// Apply the \"Added many more features\" edits to the Deadly 5 review.\n// Each pair is an exact (old, new) whole-paragraph-text match. The title\n// line appears twice (the H1 and again as a bolded run near the end) and\n// both occurrences become the same new text, so we search/replace across\n// the whole body rather than targeting individual paragraphs.\n\nconst replacements = [\n  [\n    \"Play Deadly 5 Slot Free - Review of Push Gaming's New Western Game\",\n    \"Play Deadly 5 Free - Traditional Western Slot Game\",\n  ],\n  [\n    \"Impressive graphics and attention to detail\",\n    \"Traditional Western theme with spectacular graphics and attention to detail\",\n  ],\n  [\n    \"Significant range of betting options\",\n    \"Medium volatility creates a balanced gameplay experience\",\n  ],\n  [\n    \"High RTP of 96.39%\",\n    \"Impressive 96.39% RTP offers higher average payouts\",\n  ],\n  [\n    \"Free spins with extended wild symbol and multiplier\",\n    \"Bonus features provide opportunities to win big\",\n  ],\n  [\n    \"Only 20 fixed pay lines\",\n    \"Limited betting options for players\",\n  ],\n  [\n    \"Scatter symbols are not very common\",\n    \"May not appeal to players who are not fans of the Western genre\",\n  ],\n  [\n    \"Read our review of Deadly 5 slot game, a newly released online slot game by Push Gaming. Play for free and enjoy Western-themed graphics and features.\",\n    \"Read our review of Deadly 5 and play this traditional Western slot game for free!\",\n  ],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, {\n    matchCase: true,\n    matchWholeWord: false,\n  });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Apply the \"Added many more features\" edits to the Deadly 5 review.\n# Each pair is an exact (old, new) whole-text match against a run's text in\n# the document; wdReplaceAll (2) replaces every occurrence (the title line\n# appears twice - as the H1 and again as a bolded run near the end - and\n# both must become the same new text).\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"Play Deadly 5 Slot Free - Review of Push Gaming's New Western Game\", \"Play Deadly 5 Free - Traditional Western Slot Game\"),\n    @(\"Impressive graphics and attention to detail\", \"Traditional Western theme with spectacular graphics and attention to detail\"),\n    @(\"Significant range of betting options\", \"Medium volatility creates a balanced gameplay experience\"),\n    @(\"High RTP of 96.39%\", \"Impressive 96.39% RTP offers higher average payouts\"),\n    @(\"Free spins with extended wild symbol and multiplier\", \"Bonus features provide opportunities to win big\"),\n    @(\"Only 20 fixed pay lines\", \"Limited betting options for players\"),\n    @(\"Scatter symbols are not very common\", \"May not appeal to players who are not fans of the Western genre\"),\n    @(\"Read our review of Deadly 5 slot game, a newly released online slot game by Push Gaming. Play for free and enjoy Western-themed graphics and features.\", \"Read our review of Deadly 5 and play this traditional Western slot game for free!\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute([ref]$oldText, $false, $true, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
